# Updated cryptos list on Tue Jul 30 05:17:34 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.661.24'
$ws.Range('E2').Value = '  -4.05%  '
$ws.Range('D3').Value = '3.311.78'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').Value = '183.33'
$ws.Range('E5').Value = '  -5.18%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '572.74'
$ws.Range('E6').Value = '  -3.21%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -1.41%  '
$ws.Range('E9').Value = '  -3.58%  '
$ws.Range('E10').Value = '  -1.99%  '
$ws.Range('D11').Value = '0.403'
$ws.Range('E11').Value = '  -4.23%  '
$ws.Range('D12').Value = '3.889.57'
$ws.Range('E12').Value = '  -0.93%  '
$ws.Range('D13').Value = '0.138'
$ws.Range('D14').Value = '27.17'
$ws.Range('E14').Value = '  -4.14%  '
$ws.Range('D15').Value = '66.736.09'
$ws.Range('E15').Value = '  -3.85%  '
$ws.Range('E16').Value = '  -2.50%  '
$ws.Range('D17').Value = '3.322.16'
$ws.Range('E17').Value = '  -1.01%  '
$ws.Range('D18').Value = '13.70'
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('D19').Value = '433.99'
$ws.Range('E19').Value = '  -0.95%  '
$ws.Range('D20').Value = '5.67'
$ws.Range('E20').Value = '  -2.53%  '
$ws.Range('D21').Value = '7.62'
$ws.Range('E21').Value = '  -1.92%  '
$ws.Range('D22').Value = '73.76'
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = '0.516'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  -2.97%  '
$ws.Range('D26').Value = '0.193'
$ws.Range('E26').Value = '  +0.86%  '
$ws.Range('D27').Value = '9.04'
$ws.Range('E27').Value = '  -5.36%  '
$ws.Range('E28').Value = '  -0.54%  '
$ws.Range('E29').Value = '  -1.87%  '
$ws.Range('D30').Value = '22.78'
$ws.Range('E30').Value = '  -1.24%  '
$ws.Range('D31').Value = '5.32'
$ws.Range('E31').Value = '  -4.56%  '
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = '1.24'
$ws.Range('E33').Value = '  -3.03%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').Value = '6.78'
$ws.Range('E34').Value = '  -3.39%  '
$ws.Range('E35').Value = '  -1.29%  '
$ws.Range('D36').Value = '160.47'
$ws.Range('E36').Value = '  -2.39%  '
$ws.Range('E37').Value = '  -4.21%  '
$ws.Range('D38').Value = '27.24'
$ws.Range('E38').Value = '  +0.50%  '
$ws.Range('D39').Value = '2.818.28'
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('D40').Value = '0.790'
$ws.Range('E40').Value = '  -2.29%  '
$ws.Range('D41').Value = '4.44'
$ws.Range('E41').Value = '  -2.64%  '
$ws.Range('D42').Value = '6.21'
$ws.Range('E42').Value = '  -4.33%  '
$ws.Range('E43').Value = '  -1.77%  '
$ws.Range('D44').Value = '40.18'
$ws.Range('E44').Value = '  -1.45%  '
$ws.Range('D45').Value = '24.45'
$ws.Range('E45').Value = '  -3.42%  '
$ws.Range('D46').Value = '2.35'
$ws.Range('E46').Value = '  -6.80%  '
$ws.Range('D47').Value = '320.51'
$ws.Range('E47').Value = '  -6.96%  '
$ws.Range('D48').Value = '0.0272'
$ws.Range('E48').Value = '  -3.72%  '
$ws.Range('D49').Value = '0.986'
$ws.Range('E49').Value = '  -2.25%  '
$ws.Range('D50').Value = '6.16'
$ws.Range('E50').Value = '  -1.97%  '
$ws.Range('D51').Value = '0.100'
$ws.Range('E51').Value = '  -1.03%  '
